$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1424448227877377
$ws.Range("D2").Value = 0.1538830668338065
$ws.Range("E2").Value = 0.1615681961842341
$ws.Range("F2").Value = 0.9766215052929894
$ws.Range("G2").Value = 0.002421554833943978
$ws.Range("J2").Value = 0.2335619647207068
$ws.Range("M2").Value = 5.60687900451768
$ws.Range("O2").Value = 2.47296806659557
$ws.Range("B3").Value = 0.1329302348927399
$ws.Range("D3").Value = 0.1555917809322693
$ws.Range("E3").Value = 0.1577078460882788
$ws.Range("F3").Value = 1.00199806685967
$ws.Range("G3").Value = 0.002426071402750196
$ws.Range("J3").Value = 0.2210698601333831
$ws.Range("M3").Value = 4.9278393971139
$ws.Range("O3").Value = 2.511344351416767
$ws.Range("B4").Value = 0.1271576539167114
$ws.Range("D4").Value = 0.1567560085854325
$ws.Range("E4").Value = 0.1555030991260651
$ws.Range("F4").Value = 1.019011050025611
$ws.Range("G4").Value = 0.002428976965463307
$ws.Range("J4").Value = 0.2136342219636447
$ws.Range("M4").Value = 4.508989111625823
$ws.Range("O4").Value = 2.538473015305158
$ws.Range("B5").Value = 0.1248228521123451
$ws.Range("D5").Value = 0.1572593140157252
$ws.Range("E5").Value = 0.1546458931454069
$ws.Range("F5").Value = 1.026300364996601
$ws.Range("G5").Value = 0.002430194418517894
$ws.Range("J5").Value = 0.2106622849612592
$ws.Range("M5").Value = 4.337822420834954
$ws.Range("O5").Value = 2.550415171844747
$ws.Range("B6").Value = 0.1244362248811797
$ws.Range("D6").Value = 0.157344629650467
$ws.Range("E6").Value = 0.1545060354371373
$ws.Range("F6").Value = 1.027532169282487
$ws.Range("G6").Value = 0.002430398597274582
$ws.Range("J6").Value = 0.2101722855981052
$ws.Range("M6").Value = 4.309371256298647
$ws.Range("O6").Value = 2.552451451505618
$ws.Range("B7").Value = 0.1271260946575126
$ws.Range("D7").Value = 0.1567626795086312
$ws.Range("E7").Value = 0.1554913720177922
$ws.Range("F7").Value = 1.019107918048888
$ws.Range("G7").Value = 0.002428993248955146
$ws.Range("J7").Value = 0.213593907013987
$ws.Range("M7").Value = 4.506682647289523
$ws.Range("O7").Value = 2.538630492544883
$ws.Range("B8").Value = 0.1391498514531122
$ws.Range("D8").Value = 0.1544483094521283
$ws.Range("E8").Value = 0.1602025567087537
$ws.Range("F8").Value = 0.9850720228673566
$ws.Range("G8").Value = 0.00242308475031481
$ws.Range("J8").Value = 0.2292054767660545
$ws.Range("M8").Value = 5.373143669790153
$ws.Range("O8").Value = 2.485454436807942
$ws.Range("B9").Value = 0.1632755997721205
$ws.Range("D9").Value = 0.150826009249684
$ws.Range("E9").Value = 0.1707736206860346
$ws.Range("F9").Value = 0.9298492921190089
$ws.Range("G9").Value = 0.002412542586825175
$ws.Range("J9").Value = 0.2617234804396418
$ws.Range("M9").Value = 7.057192597014875
$ws.Range("O9").Value = 2.409897481932347
$ws.Range("B10").Value = 0.1813318817857237
$ws.Range("D10").Value = 0.1487280798060482
$ws.Range("E10").Value = 0.1793820175676686
$ws.Range("F10").Value = 0.8965233744792727
$ws.Range("G10").Value = 0.002405425468884395
$ws.Range("J10").Value = 0.2868410619461912
$ws.Range("M10").Value = 8.285668969970743
$ws.Range("O10").Value = 2.372491108661535
$ws.Range("B11").Value = 0.1896177386592512
$ws.Range("D11").Value = 0.1478971968063263
$ws.Range("E11").Value = 0.183488063403324
$ws.Range("F11").Value = 0.882984249735749
$ws.Range("G11").Value = 0.002402322265558181
$ws.Range("J11").Value = 0.2985498045417785
$ws.Range("M11").Value = 8.84274510945886
$ws.Range("O11").Value = 2.359534969321601
$ws.Range("B12").Value = 0.1927656611223512
$ws.Range("D12").Value = 0.1476004258944812
$ws.Range("E12").Value = 0.1850708444095375
$ws.Range("F12").Value = 0.8780944770312615
$ws.Range("G12").Value = 0.002401166346932569
$ws.Range("J12").Value = 0.3030255708276286
$ws.Range("M12").Value = 9.053450656428538
$ws.Range("O12").Value = 2.355223417806769
$ws.Range("B13").Value = 0.1920872454431333
$ws.Range("D13").Value = 0.1476635443611656
$ws.Range("E13").Value = 0.1847287132261783
$ws.Range("F13").Value = 0.8791369617478466
$ws.Range("G13").Value = 0.002401414442876239
$ws.Range("J13").Value = 0.3020597489515637
$ws.Range("M13").Value = 9.008082317804622
$ws.Range("O13").Value = 2.35612536673753
$ws.Range("B14").Value = 0.1898765152999005
$ws.Range("D14").Value = 0.1478724225541725
$ws.Range("E14").Value = 0.1836177169223632
$ws.Range("F14").Value = 0.8825771877054009
$ws.Range("G14").Value = 0.002402226783446077
$ws.Range("J14").Value = 0.2989171808151525
$ws.Range("M14").Value = 8.86008489141409
$ws.Range("O14").Value = 2.359168275280723
$ws.Range("B15").Value = 0.1885237109988509
$ws.Range("D15").Value = 0.1480026966300159
$ws.Range("E15").Value = 0.1829408522572038
$ws.Range("F15").Value = 0.8847154383446991
$ws.Range("G15").Value = 0.00240272686138177
$ws.Range("J15").Value = 0.2969977654434786
$ws.Range("M15").Value = 8.769400232868691
$ws.Range("O15").Value = 2.361109909134598
$ws.Range("B16").Value = 0.1807918209862436
$ws.Range("D16").Value = 0.1487848739149769
$ws.Range("E16").Value = 0.1791175509937091
$ws.Range("F16").Value = 0.8974411565003635
$ws.Range("G16").Value = 0.002405630964131436
$ws.Range("J16").Value = 0.2860816628537037
$ws.Range("M16").Value = 8.249227768664127
$ws.Range("O16").Value = 2.373420445626437
$ws.Range("B17").Value = 0.1760669158413037
$ws.Range("D17").Value = 0.1492964159502677
$ws.Range("E17").Value = 0.1768211578577592
$ws.Range("F17").Value = 0.9056660551924125
$ws.Range("G17").Value = 0.002407446869872795
$ws.Range("J17").Value = 0.2794582614743746
$ws.Range("M17").Value = 7.929671259609279
$ws.Range("O17").Value = 2.382020129966463
$ws.Range("B18").Value = 0.1733560546019959
$ws.Range("D18").Value = 0.1496022564410211
$ws.Range("E18").Value = 0.1755181840630868
$ws.Range("F18").Value = 0.9105491281243374
$ws.Range("G18").Value = 0.002408503989944953
$ws.Range("J18").Value = 0.2756751775657875
$ws.Range("M18").Value = 7.745704301134651
$ws.Range("O18").Value = 2.387347954499518
$ws.Range("B19").Value = 0.1724393710299239
$ws.Range("D19").Value = 0.1497078004079881
$ws.Range("E19").Value = 0.175080068212722
$ws.Range("F19").Value = 0.9122284909979896
$ws.Range("G19").Value = 0.002408864090982412
$ws.Range("J19").Value = 0.2743988076873762
$ws.Range("M19").Value = 7.683387475022244
$ws.Range("O19").Value = 2.38921706219341
$ws.Range("B20").Value = 0.1765691891744297
$ws.Range("D20").Value = 0.1492407585362443
$ws.Range("E20").Value = 0.1770637612258241
$ws.Range("F20").Value = 0.9047747056681317
$ws.Range("G20").Value = 0.002407252254508546
$ws.Range("J20").Value = 0.2801605799740798
$ws.Range("M20").Value = 7.963705785129434
$ws.Range("O20").Value = 2.381065115790932
$ws.Range("B21").Value = 0.1905255829701389
$ws.Range("D21").Value = 0.147810584155863
$ws.Range("E21").Value = 0.1839432811699666
$ws.Range("F21").Value = 0.8815602395023419
$ws.Range("G21").Value = 0.002401987659268212
$ws.Range("J21").Value = 0.2998390810163869
$ws.Range("M21").Value = 8.903561970232431
$ws.Range("O21").Value = 2.358258273775277
$ws.Range("B22").Value = 0.1997065850139137
$ws.Range("D22").Value = 0.1469800784874025
$ws.Range("E22").Value = 0.1886024239513517
$ws.Range("F22").Value = 0.8677727619983528
$ws.Range("G22").Value = 0.002398658776785032
$ws.Range("J22").Value = 0.3129451041353377
$ws.Range("M22").Value = 9.516376056101763
$ws.Range("O22").Value = 2.34682389860842
$ws.Range("B23").Value = 0.1948010814095937
$ws.Range("D23").Value = 0.1474137625359191
$ws.Range("E23").Value = 0.186100637808714
$ws.Range("F23").Value = 0.8750033480019397
$ws.Range("G23").Value = 0.002400425274694197
$ws.Range("J23").Value = 0.3059273153248228
$ws.Range("M23").Value = 9.189434545315464
$ws.Range("O23").Value = 2.35260542349647
$ws.Range("B24").Value = 0.1763420941068006
$ws.Range("D24").Value = 0.1492658846682886
$ws.Range("E24").Value = 0.1769540265468521
$ws.Range("F24").Value = 0.9051772039613439
$ws.Range("G24").Value = 0.0024073401991857
$ws.Range("J24").Value = 0.2798429846364456
$ws.Range("M24").Value = 7.948319554267357
$ws.Range("O24").Value = 2.381495683216144
$ws.Range("B25").Value = 0.1566906398686001
$ws.Range("D25").Value = 0.1517073783763863
$ws.Range("E25").Value = 0.1677682861673588
$ws.Range("F25").Value = 0.9435321719888705
$ws.Range("G25").Value = 0.002415283583991111
$ws.Range("J25").Value = 0.2527161236283035
$ws.Range("M25").Value = 6.603190885150184
$ws.Range("O25").Value = 2.427201531225904
